# The presentation currently ships two theme parts:
#   - the theme behind the (single) slide master uses the "Integral" color
#     scheme (dk2=455F51, lt2=E3DED1, accent1=99CB38, ...)
#   - the theme behind the notes master uses the default "Office Theme"
#     color scheme (dk2=44546A, lt2=E7E6E6, accent1=5B9BD5, ...)
#
# The target edit swaps the two color schemes: the slide master's theme
# should end up carrying the "Office Theme" palette (and vice versa for the
# notes master's theme). The font scheme and format scheme of both themes
# are already identical, so only the 12 theme colors need to change.
#
# Apply the swap through the Design/Theme COM surface: walk to the active
# presentation's slide master theme and overwrite each of its 12
# ThemeColorScheme entries (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# with the "Office Theme" values.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

function ConvertTo-ComColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Index order exposed by ThemeColorScheme.Item(1..12):
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-ComColor($officeThemeColors[$i - 1])
}
